$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 2 to make room for the FAPs and sCs rows
$ws.Rows("3:4").Insert()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Rhbdl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1533166666666667
$ws.Range("H2").Value = 0.45995
$ws.Range("I2").Value = 0.1117088182569538
$ws.Range("J2").Value = 0.1117088182569538
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.825757333333333
$ws.Range("N2").Value = 5.477272
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.2799190284888889
$ws.Range("R2").Value = 2.5192712564
$ws.Range("S2").Value = 0.1117088182569538
$ws.Range("T2").Value = 0.1117088182569538

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Rhbdl2"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9559960000000002
$ws.Range("H3").Value = 2.867988
$ws.Range("I3").Value = 0.6965529954454279
$ws.Range("J3").Value = 0.6965529954454278
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.825757333333333
$ws.Range("N3").Value = 5.477272
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.745416707637334
$ws.Range("R3").Value = 15.708750368736
$ws.Range("S3").Value = 0.6965529954454279
$ws.Range("T3").Value = 0.6965529954454278

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Rhbdl2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2631543333333333
$ws.Range("H4").Value = 0.789463
$ws.Range("I4").Value = 0.1917381862976183
$ws.Range("J4").Value = 0.1917381862976183
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.825757333333333
$ws.Range("N4").Value = 5.477272
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.4804559538817778
$ws.Range("R4").Value = 4.324103584936
$ws.Range("S4").Value = 0.1917381862976183
$ws.Range("T4").Value = 0.1917381862976183
